# Auto-generated edit script: updates market-price-derived cells (columns H-N)
# on sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR to reflect a refreshed
# data snapshot from the scheduled Universalis price-sync runner.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 347.83334
$ws.Range("J4").Value = 375
$ws.Range("L4").Value = 375
$ws.Range("N4").Value = -603
$ws.Range("H9").Value = 10000118
$ws.Range("I9").Value = 16666741
$ws.Range("K9").Value = 16666741
$ws.Range("M9").Value = -16666572
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H38").Value = 535.8
$ws.Range("I38").Value = 317.55554
$ws.Range("J38").Value = 2500
$ws.Range("K38").Value = 952.66662
$ws.Range("L38").Value = 7500
$ws.Range("M38").Value = -580.66662
$ws.Range("N38").Value = -8244
$ws.Range("H40").Value = 1833.3334
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H43").Value = 1230
$ws.Range("I43").Value = 1434
$ws.Range("K43").Value = 1434
$ws.Range("M43").Value = -1365
$ws.Range("H58").Value = 73117.14
$ws.Range("J58").Value = 102170.5
$ws.Range("L58").Value = 306511.5
$ws.Range("N58").Value = -306811.5
$ws.Range("H87").Value = 27488.88
$ws.Range("J87").Value = 27488.88
$ws.Range("L87").Value = 27488.88
$ws.Range("N87").Value = -29984.88
$ws.Range("H90").Value = 27488.88
$ws.Range("J90").Value = 27488.88
$ws.Range("L90").Value = 82466.64
$ws.Range("N90").Value = -94946.64
$ws.Range("H112").Value = 1412.4142
$ws.Range("J112").Value = 1459.1398
$ws.Range("L112").Value = 4377.4194
$ws.Range("N112").Value = -6593.4194
$ws.Range("H129").Value = 5682863
$ws.Range("J129").Value = 982.125
$ws.Range("L129").Value = 2946.375
$ws.Range("N129").Value = -12946.375
$ws.Range("H132").Value = 3574033.2
$ws.Range("I132").Value = 4002012
$ws.Range("K132").Value = 12006036
$ws.Range("M132").Value = -12003506
$ws.Range("H138").Value = 5391.89
$ws.Range("I138").Value = 2973.16
$ws.Range("J138").Value = 6198.1333
$ws.Range("K138").Value = 8919.48
$ws.Range("L138").Value = 18594.3999
$ws.Range("M138").Value = -3779.48
$ws.Range("N138").Value = -28874.3999
$ws.Range("H141").Value = 516416.75
$ws.Range("I141").Value = 1770.7222
$ws.Range("J141").Value = 2060354.9
$ws.Range("K141").Value = 5312.1666
$ws.Range("L141").Value = 6181064.699999999
$ws.Range("M141").Value = -132.1665999999996
$ws.Range("N141").Value = -6191424.699999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 14330002
$ws.Range("I6").Value = 50005000
$ws.Range("J6").Value = 60001.8
$ws.Range("K6").Value = 50005000
$ws.Range("L6").Value = 60001.8
$ws.Range("M6").Value = -50004827
$ws.Range("N6").Value = -60347.8
$ws.Range("H32").Value = 19030.291
$ws.Range("I32").Value = 15997.56
$ws.Range("J32").Value = 31666.666
$ws.Range("K32").Value = 15997.56
$ws.Range("L32").Value = 31666.666
$ws.Range("M32").Value = -15710.56
$ws.Range("N32").Value = -32240.666
$ws.Range("H133").Value = 28800
$ws.Range("J133").Value = 28800
$ws.Range("L133").Value = 28800
$ws.Range("N133").Value = -33860

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H141").Value = 30000
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2861162.2
$ws.Range("I31").Value = 4002395.5
$ws.Range("K31").Value = 4002395.5
$ws.Range("M31").Value = -4002100.5
$ws.Range("H34").Value = 2861162.2
$ws.Range("I34").Value = 4002395.5
$ws.Range("K34").Value = 4002395.5
$ws.Range("M34").Value = -4002193.5
$ws.Range("H134").Value = 2648.9512
$ws.Range("I134").Value = 2089.08
$ws.Range("J134").Value = 3523.75
$ws.Range("K134").Value = 6267.24
$ws.Range("L134").Value = 10571.25
$ws.Range("M134").Value = -3732.24
$ws.Range("N134").Value = -15641.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 166.66667
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 300
$ws.Range("L7").Value = 600
$ws.Range("M7").Value = -188
$ws.Range("N7").Value = -824
$ws.Range("H34").Value = 6530.722
$ws.Range("I34").Value = 114.28571
$ws.Range("J34").Value = 10613.909
$ws.Range("K34").Value = 342.85713
$ws.Range("L34").Value = 31841.727
$ws.Range("M34").Value = -258.85713
$ws.Range("N34").Value = -32009.727
$ws.Range("H39").Value = 1875
$ws.Range("H55").Value = 1472.2222
$ws.Range("I55").Value = 670
$ws.Range("J55").Value = 2475
$ws.Range("K55").Value = 2010
$ws.Range("L55").Value = 7425
$ws.Range("M55").Value = -1833
$ws.Range("N55").Value = -7779
$ws.Range("H114").Value = 1361.3077
$ws.Range("I114").Value = 775
$ws.Range("J114").Value = 1621.8889
$ws.Range("K114").Value = 2325
$ws.Range("L114").Value = 4865.6667
$ws.Range("M114").Value = 929
$ws.Range("N114").Value = -11373.6667
$ws.Range("H129").Value = 18710.2
$ws.Range("I129").Value = 2341.7144
$ws.Range("J129").Value = 33032.625
$ws.Range("K129").Value = 7025.1432
$ws.Range("L129").Value = 99097.875
$ws.Range("M129").Value = -2025.1432
$ws.Range("N129").Value = -109097.875
$ws.Range("H131").Value = 1512.6842
$ws.Range("J131").Value = 1176.3518
$ws.Range("L131").Value = 3529.0554
$ws.Range("N131").Value = -13609.0554

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6160.091
$ws.Range("I132").Value = 8120.364
$ws.Range("J132").Value = 4199.8184
$ws.Range("K132").Value = 24361.092
$ws.Range("L132").Value = 12599.4552
$ws.Range("M132").Value = -21831.092
$ws.Range("N132").Value = -17659.4552

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1581.7273
$ws.Range("I22").Value = 520.2
$ws.Range("J22").Value = 2466.3333
$ws.Range("K22").Value = 520.2
$ws.Range("L22").Value = 2466.3333
$ws.Range("M22").Value = -225.2
$ws.Range("N22").Value = -3056.3333
$ws.Range("H27").Value = 1581.7273
$ws.Range("I27").Value = 520.2
$ws.Range("J27").Value = 2466.3333
$ws.Range("K27").Value = 520.2
$ws.Range("L27").Value = 2466.3333
$ws.Range("M27").Value = -413.2
$ws.Range("N27").Value = -2680.3333

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 28565
$ws.Range("J133").Value = 28565
$ws.Range("L133").Value = 28565
$ws.Range("N133").Value = -38685
